$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Bloc section: merge the two runs "I" + "d" into a single run "Id".
#    (paragraph 25 in the original layout)
# ---------------------------------------------------------------------------
$pId = $d.Paragraphs.Item(25)
if ($pId.Range.Text.Trim() -eq "Id") {
    $pId.Range.Text = "Id"
}

# ---------------------------------------------------------------------------
# 2) Room section: remove the "Events {}" bullet - it is superseded by the
#    new "events = {}" variable declared under the World chapter (step 4).
#    (paragraph 21 in the original layout)
# ---------------------------------------------------------------------------
$pEvents = $d.Paragraphs.Item(21)
if ($pEvents.Range.Text.Trim() -eq "Events {}") {
    $pEvents.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Room section: the "_GoBack" bookmark now wraps the "X" bullet instead of
#    the "Y" bullet (paragraphs 16 and 17 in the original layout).
# ---------------------------------------------------------------------------
$pX = $d.Paragraphs.Item(16)
$pY = $d.Paragraphs.Item(17)
if ($pX.Range.Text.Trim() -eq "X" -and $pY.Range.Text.Trim() -eq "Y") {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
    $xWordRange = $pX.Range.Words.Item(1)
    $d.Bookmarks.Add("_GoBack", $xWordRange)
}

# ---------------------------------------------------------------------------
# 4) World section: add a new "events = {}" bullet right after "rooms = {}",
#    mirroring the existing room dictionary declaration.
# ---------------------------------------------------------------------------
$pRooms = $d.Paragraphs.Item(8)
if ($pRooms.Range.Text.Trim() -eq "rooms = {}") {
    $pRooms.Range.InsertParagraphAfter()
    $pNew = $d.Paragraphs.Item(9)
    $pNew.Range.Text = "events = {}"
}
